$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The boolean flag previously stored in K2 (is_alumni) is being moved to
# L2 (is_student): clear K2 and set L2 = TRUE, matching the "edit profile
# icon" fix described in the commit message.
$ws.Range("K2").ClearContents()
$ws.Range("L2").Value = $true

# Update the active selection to the new cell, as captured in the diff.
$ws.Range("L2").Select()
